# Update the captured run timestamp embedded in the OLS regression summary
# text on each of the three sheets ("5", "4", "3").  The statsmodels output
# was re-generated the next day, so every occurrence of the old Date:/Time:
# stamp is replaced with the new one; everything else in the summary block
# (coefficients, AIC/BIC, etc.) stays untouched.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $cell = $ws.Range("B2")
    $text = $cell.Value2

    if ($text -ne $null -and $text.Contains("Date:")) {
        $text = $text.Replace("Date:                Sat, 28 Dec 2019", "Date:                Sun, 29 Dec 2019")
        $text = $text.Replace("Time:                        20:59:42", "Time:                        16:11:11")
        $text = $text.Replace("Time:                        20:59:43", "Time:                        16:11:11")
        $cell.Value2 = $text
    }
}
